# Learning Journal 1 — update the Journal URL hyperlink text to show the
# full GitHub URL, and re-flow the two table-cell paragraphs whose text
# wraps across a rendered page break (the words move from one side of
# the <w:lastRenderedPageBreak/> marker to the other; the combined text
# of each cell is unchanged).

$d = $word.ActiveDocument

# --- 1. Hyperlink display text -------------------------------------------------
# The hyperlink currently displays "Github URL for the file." (split across
# two runs with a spell-check proofErr wrapper). Replace the *displayed*
# text with the full target URL, collapsing it to a single run.
$hyperlink = $d.Hyperlinks.Item(1)
$hyperlink.TextToDisplay = "https://github.com/tarekFerdous/softwareProjectManagementSOEN6841/blob/f6e400fd6144c512f8c9a877997be4c5aaaacd25/Learning%20Journals/Learning%20Journal%201.pdf"

# --- 2. Table cell re-flow -------------------------------------------------
$table = $d.Tables.Item(1)

# Cell (row 2, col 1): "... COCOMO. These |methods help plan resources ..."
# The page-break now falls right after "These " instead of after "plan ".
$cell1 = $table.Cell(2, 1)
$cell1Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4522EB27" w14:textId="57F45385" w:rsidR="005B7241" w:rsidRDefault="001E61FB"><w:r w:rsidRPr="001E61FB"><w:t xml:space="preserve">I learned how software projects are initiated, including creating project charters, defining the scope, and setting SMART objectives (Specific, Measurable, Achievable, Relevant, and Time-bound). I also explored techniques for estimating project effort and costs, like Function Point Analysis (FPA) and COCOMO. These </w:t></w:r><w:r w:rsidRPr="001E61FB"><w:lastRenderedPageBreak/><w:t>methods help plan resources and budgets based on the size and complexity of a software project.</w:t></w:r></w:p>'
$cell1.Range.InsertXML($cell1Xml)

# Cell (row 2, col 2): both runs keep their own <w:lastRenderedPageBreak/>;
# only the text boundary between them shifts from "... as |the project..."
# to "... adjusting |estimates as the project progresses."
$cell2 = $table.Cell(2, 2)
$cell2Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="67A41AE5" w14:textId="1A85BF71" w:rsidR="005B7241" w:rsidRDefault="001E61FB"><w:r w:rsidRPr="001E61FB"><w:lastRenderedPageBreak/><w:t xml:space="preserve">In real-world projects, I’ll use project charters to clearly define objectives and scope. Effort estimation methods like COCOMO and FPA will help me calculate the time and resources needed for development. These techniques are especially helpful in Agile projects, where requirements often change, as they provide flexibility in adjusting </w:t></w:r><w:r w:rsidRPr="001E61FB"><w:lastRenderedPageBreak/><w:t>estimates as the project progresses.</w:t></w:r></w:p>'
$cell2.Range.InsertXML($cell2Xml)

Write-Output "Applied hyperlink text + table re-flow edits."
